$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Last Updated" timestamp
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 10:42 AM"

# ---------------------------------------------------------------------------
# 2. Stock List sheet: the screener refreshed - the top two rows (MIDWESTLTD,
#    CAPTRU-RE1) dropped off the list, every remaining row shifted up by two
#    positions, and two new rows (SMARTWORKS, TRAVELFOOD) were appended at
#    the bottom.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Stock List")

# Remove the two rows that fell out of the list; everything below shifts up.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Append the two new rows that now appear at the bottom of the list.
$icon = "📋"

$ws.Range("A75").Value = $icon
$ws.Range("B75").Value = "SMARTWORKS"
$ws.Range("C75").Value = "SMARTWORKS"
$ws.Range("D75").Value = 606.65
$ws.Range("E75").Value = 2.0867
$ws.Range("F75").Value = "N/A"
$ws.Range("G75").Value = "N/A"
$ws.Range("H75").Value = 6931.2448

$ws.Range("A76").Value = $icon
$ws.Range("B76").Value = "TRAVELFOOD"
$ws.Range("C76").Value = "TRAVELFOOD"
$ws.Range("D76").Value = 1316.3
$ws.Range("E76").Value = 0.1141
$ws.Range("F76").Value = "N/A"
$ws.Range("G76").Value = "N/A"
$ws.Range("H76").Value = 17332.9705
